# Update countries & provincias Spain
# Applies latest COVID-19 snapshot figures to the "Pais" sheet.
# - Kuwait's case counts overtook Panama/Noruega/Chequia, so Kuwait now
#   sits right after Egipto; Panama, Noruega and Chequia each shift down
#   one row. Australia (which stays in place) also gets refreshed figures.
# - Benin's case counts overtook Sierra Leona/Vietnam, so Benin now sits
#   right after Republica del Chad; Sierra Leona and Vietnam each shift
#   down one row.
# - A handful of other rows (Banglades, Uzbekistan, Malta) simply get
#   updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Banglades (row 39): refreshed "Recuperados" / "Casos criticos" ---
$ws.Range("D39").Value = 2650
$ws.Range("E39").Value = 11779

# --- Kuwait overtakes Panama / Noruega / Chequia (rows 49-52) ---
# Row 49 becomes Kuwait with brand-new totals.
$ws.Range("A49").Value = "Kuwait"
$ws.Range("B49").Value = 8688
$ws.Range("C49").Value = 1065
$ws.Range("D49").Value = 2729
$ws.Range("E49").Value = 5901
$ws.Range("F49").Value = 114
$ws.Range("G49").Value = 9
$ws.Range("H49").Value = 58

# Row 50 becomes Panama (previously row 49's figures).
$ws.Range("A50").Value = "Panama"
$ws.Range("B50").Value = 8282
$ws.Range("C50").Value = 212
$ws.Range("D50").Value = 4501
$ws.Range("E50").Value = 3544
$ws.Range("F50").Value = 85
$ws.Range("G50").Value = 6
$ws.Range("H50").Value = 237

# Row 51 becomes Noruega (previously row 50's figures).
$ws.Range("A51").Value = "Noruega"
$ws.Range("B51").Value = 8099
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 32
$ws.Range("E51").Value = 7848
$ws.Range("F51").Value = 24
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 219

# Row 52 becomes Chequia (previously row 51's figures).
$ws.Range("A52").Value = "Chequia"
$ws.Range("B52").Value = 8095
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 4448
$ws.Range("E52").Value = 3371
$ws.Range("F52").Value = 40
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 276

# Row 53 stays Australia, but gets refreshed totals.
$ws.Range("B53").Value = 6941
$ws.Range("C53").Value = 12
$ws.Range("D53").Value = 6163
$ws.Range("E53").Value = 681

# Row 56 stays Finlandia, refreshed totals.
$ws.Range("B56").Value = 5962
$ws.Range("C56").Value = 82
$ws.Range("E56").Value = 1695

# --- Uzbekistan (row 75): refreshed "Recuperados" / "Casos criticos" ---
$ws.Range("D75").Value = 1870
$ws.Range("E75").Value = 531

# --- Malta (row 124): refreshed totals ---
$ws.Range("B124").Value = 496
$ws.Range("C124").Value = 6
$ws.Range("D124").Value = 433

# --- Benin overtakes Sierra Leona / Vietnam (rows 134-136) ---
# Row 134 becomes Benin with brand-new totals.
$ws.Range("A134").Value = "Benin"
$ws.Range("B134").Value = 319
$ws.Range("C134").Value = 35
$ws.Range("D134").Value = 62
$ws.Range("E134").Value = 255
$ws.Range("H134").Value = 2

# Row 135 becomes Sierra Leona (previously row 134's figures).
$ws.Range("A135").Value = "Sierra Leona"
$ws.Range("B135").Value = 291
$ws.Range("D135").Value = 58
$ws.Range("E135").Value = 215
$ws.Range("F135").Value = 0
$ws.Range("H135").Value = 18

# Row 136 becomes Vietnam (previously row 135's figures).
$ws.Range("A136").Value = "Vietnam"
$ws.Range("B136").Value = 288
$ws.Range("D136").Value = 241
$ws.Range("E136").Value = 47
$ws.Range("F136").Value = 8
$ws.Range("H136").Value = 0
